$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing header cell (G1) onto the new
# header cell (H1) so the new "Save" column header matches the style
# used by the other header cells.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 1
